# Applies the "Updated cryptos list" price/volume refresh to Sheet1 (rows 2-51).
# A leading apostrophe is used for Price (column D) values that look like plain
# numbers (e.g. "21.60", "0.0000140") so Excel stores them as text - exactly like
# the original inline-string cells - instead of silently normalising them into
# numbers (which would drop trailing zeros / switch to scientific notation).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "60.521.70"
$ws.Range("E2").Value = "  -0.22%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "2.630.89"
$ws.Range("E3").Value = "  -0.37%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.13%  "

# Row 5: BNB
$ws.Range("D5").Value = "'511.05"
$ws.Range("E5").Value = "  -0.33%  "

# Row 6: Solana
$ws.Range("D6").Value = "'154.43"
$ws.Range("E6").Value = "  -2.66%  "

# Row 7: USDC
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.60%  "

# Row 8: XRP
$ws.Range("E8").Value = "  -2.77%  "

# Row 9: LidoStakedEther
$ws.Range("D9").Value = "2.627.51"
$ws.Range("E9").Value = "  -1.98%  "

# Row 10: Toncoin
$ws.Range("D10").Value = "'6.71"
$ws.Range("E10").Value = "  +3.28%  "

# Row 11: Dogecoin
$ws.Range("E11").Value = "  -1.01%  "

# Row 12: Cardano
$ws.Range("D12").Value = "'0.346"
$ws.Range("E12").Value = "  -0.66%  "

# Row 13: TRON
$ws.Range("E13").Value = "  +1.51%  "

# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.083.16"
$ws.Range("E14").Value = "  -1.02%  "

# Row 15: WrappedBTC
$ws.Range("D15").Value = "60.487.34"
$ws.Range("E15").Value = "  -0.31%  "

# Row 16: Avalanche
$ws.Range("D16").Value = "'21.60"
$ws.Range("E16").Value = "  -1.64%  "

# Row 17: ShibaInu
$ws.Range("D17").Value = "'0.0000140"
$ws.Range("E17").Value = "  -0.59%  "

# Row 18: WrappedEther
$ws.Range("D18").Value = "2.612.44"
$ws.Range("E18").Value = "  -2.45%  "

# Row 19: Polkadot
$ws.Range("E19").Value = "  -1.09%  "

# Row 20: BitcoinCash
$ws.Range("D20").Value = "'350.48"
$ws.Range("E20").Value = "  +0.37%  "

# Row 21: Chainlink
$ws.Range("D21").Value = "'10.60"
$ws.Range("E21").Value = "  +0.22%  "

# Row 22: Uniswap
$ws.Range("E22").Value = "  -1.02%  "

# Row 23: Dai
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  +0.14%  "

# Row 24: Litecoin
$ws.Range("D24").Value = "'60.60"
$ws.Range("E24").Value = "  +0.38%  "

# Row 25: Polygon
$ws.Range("E25").Value = "  -0.59%  "

# Row 26: Kaspa
$ws.Range("E26").Value = "  -1.27%  "

# Row 27: Binance-PegBSC-USD
$ws.Range("D27").Value = "'0.997"
$ws.Range("E27").Value = "  +0.47%  "

# Row 28: PEPE
$ws.Range("D28").Value = "0.0₃0842"
$ws.Range("E28").Value = "  -4.14%  "

# Row 29: InternetComputer(DFINITY)
$ws.Range("D29").Value = "'7.37"
$ws.Range("E29").Value = "  -2.97%  "

# Row 30: USDe
$ws.Range("E30").Value = "  +0.36%  "

# Row 31: EthereumClassic
$ws.Range("D31").Value = "'19.44"
$ws.Range("E31").Value = "  -0.93%  "

# Row 32: PancakeSwap
$ws.Range("E32").Value = "  -0.70%  "

# Row 33: Monero
$ws.Range("D33").Value = "'150.44"
$ws.Range("E33").Value = "  -4.48%  "

# Row 34: Aptos
$ws.Range("D34").Value = "'5.80"
$ws.Range("E34").Value = "  +0.25%  "

# Row 35: NEARProtocol
$ws.Range("D35").Value = "'4.00"
$ws.Range("E35").Value = "  -2.48%  "

# Row 36: ImmutableX
$ws.Range("E36").Value = "  -2.93%  "

# Row 37: SuiNetwork
$ws.Range("D37").Value = "'0.883"
$ws.Range("E37").Value = "  +4.23%  "

# Row 38: Stacks
$ws.Range("D38").Value = "'1.48"
$ws.Range("E38").Value = "  -1.53%  "

# Row 39: Fetch.AI
$ws.Range("D39").Value = "'0.847"
$ws.Range("E39").Value = "  -1.65%  "

# Row 40: OKB
$ws.Range("E40").Value = "  +2.81%  "

# Row 41: Filecoin
$ws.Range("D41").Value = "'3.76"
$ws.Range("E41").Value = "  -0.56%  "

# Row 42: Bittensor
$ws.Range("D42").Value = "'293.41"
$ws.Range("E42").Value = "  -6.38%  "

# Row 43: Mantle
$ws.Range("D43").Value = "'0.625"
$ws.Range("E43").Value = "  -3.35%  "

# Row 44: Stellar
$ws.Range("E44").Value = "  -0.31%  "

# Row 45: FirstDigitalUSD
$ws.Range("D45").Value = "'0.998"
$ws.Range("E45").Value = "  +0.68%  "

# Row 46: EnergySwap (was Hedera)
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'19.93"
$ws.Range("E46").Value = "  -1.43%  "

# Row 47: Hedera (was EnergySwap)
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "'0.0555"
$ws.Range("E47").Value = "  -5.09%  "

# Row 48: RenderToken
$ws.Range("D48").Value = "'4.81"
$ws.Range("E48").Value = "  -3.07%  "

# Row 49: VeChain
$ws.Range("E49").Value = "  -1.34%  "

# Row 50: WhiteBITCoin
$ws.Range("E50").Value = "  +0.06%  "

# Row 51: Maker
$ws.Range("D51").Value = "2.002.49"
$ws.Range("E51").Value = "  -3.75%  "
